$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (persistence, horizon 1)
$ws.Range("C2").Value = 0.02045644479963466
$ws.Range("D2").Value = 0.0008532684964036991
$ws.Range("E2").Value = 0.029210759942249
$ws.Range("F2").Value = 0.9832547184925169

# Row 3 (random_forest, horizon 1)
$ws.Range("C3").Value = 0.005374860246620195
$ws.Range("D3").Value = 0.0001076096118544995
$ws.Range("E3").Value = 0.010373505282907
$ws.Range("F3").Value = 0.9978856397616686

# Row 4 (neural_network, horizon 1)
$ws.Range("C4").Value = 0.01203881740834067
$ws.Range("D4").Value = 0.0002962133632098551
$ws.Range("E4").Value = 0.01721085015941557
$ws.Range("F4").Value = 0.9941865463468429
